$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename worksheet tab: "t" -> "T"
$ws.Name = "T"

# --- "Zorg"/mentor-uur section: unify Mentoruur / mentoruur / Mentoring -> mentorles ---
$ws.Range("I19").Value = "mentorles"
$ws.Range("T19").Value = "mentorles"
$ws.Range("AG20").Value = "mentorles"
$ws.Range("Y25").Value = "mentorles"

# Mentor hours reduced from 1 to 0.5 for both columns
$ws.Range("J19").Value = 0.5
$ws.Range("K19").Value = 0.5

# Capitalize the "mentorles" label on the brugklas table (becomes its own distinct string "Mentorles")
$ws.Range("A20").Value = "Mentorles"

# --- Turn static totals into live SUM formulas ---
$ws.Range("J21").Formula = "=SUM(J3:J20)"
$ws.Range("K21").Formula = "=SUM(K3:K20)"
$ws.Range("L21").Formula = "=SUM(L3:L20)"
$ws.Range("M21").Formula = "=SUM(M3:M20)"
# Unify M21's right border weight with the rest of the row (J21:L21)
$ws.Range("M21").Borders.Item(10).Weight = 2

$ws.Range("AH28").Formula = "=SUM(AH4:AH27)"
$ws.Range("AJ28").Formula = "=SUM(AJ4:AJ27)"
$ws.Range("AL28").Formula = "=SUM(AL4:AL27)"
$ws.Range("AN28").Formula = "=SUM(AN4:AN27)"
$ws.Range("AP28").Formula = "=SUM(AP4:AP27)"
$ws.Range("AR28").Formula = "=SUM(AR4:AR27)"

# Move/restore the active selection
$ws.Range("A21").Select()

$wb.Save()
